# Generate Report for Handback
#
# This localization-status report is refreshed once the de-de / zh-cn
# handback packages have been produced: the "Status" column moves from
# "In Translation" to "Handed back: in sync with en-US", the per-locale
# "Latest Target File" / "Latest Handback File" / "Latest Handback
# DateTime" columns get populated for both tracked source files, and the
# three worksheets (Overview, zh-cn, de-de) are widened so the longer
# strings remain fully visible.

$wb  = $excel.ActiveWorkbook
$ovw = $wb.Worksheets.Item("Overview")
$zh  = $wb.Worksheets.Item("zh-cn")
$de  = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

$file1Md  = "44b9f047-174a-40d9-bd14-8fad154549c5.md"
$file2Md  = "4c41ca5f-f632-419d-ae98-52b829462cdc.md"

$file1ZhXlf = "44b9f047-174a-40d9-bd14-8fad154549c5.e9159b59ecca7399cce83702d77d249501361a41.zh-cn.xlf"
$file2ZhXlf = "4c41ca5f-f632-419d-ae98-52b829462cdc.eec83b46a069f4b9888610970020d5944ea8db17.zh-cn.xlf"
$file1DeXlf = "44b9f047-174a-40d9-bd14-8fad154549c5.e9159b59ecca7399cce83702d77d249501361a41.de-de.xlf"
$file2DeXlf = "4c41ca5f-f632-419d-ae98-52b829462cdc.eec83b46a069f4b9888610970020d5944ea8db17.de-de.xlf"

$file1Url = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c786cc5c5119b299e17aeb59fc4b2e8d5ca3802d/e2e/44b9f047-174a-40d9-bd14-8fad154549c5.md"
$file2Url = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c786cc5c5119b299e17aeb59fc4b2e8d5ca3802d/e2e/4c41ca5f-f632-419d-ae98-52b829462cdc.md"

$zhHandbackTime = "2016-09-02 00:33:16"
$deHandbackTime = "2016-09-02 00:33:23"

# ---------------------------------------------------------------------------
# 1. Status: "In Translation" -> "Handed back: in sync with en-US"
#    (Overview mirrors the same shared text in its zh-cn / de-de columns.)
# ---------------------------------------------------------------------------
$ovw.Range("E2").Value = $newStatus
$ovw.Range("F2").Value = $newStatus
$ovw.Range("E3").Value = $newStatus
$ovw.Range("F3").Value = $newStatus

$zh.Range("C2").Value = $newStatus
$zh.Range("C3").Value = $newStatus

$de.Range("C2").Value = $newStatus
$de.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------------
# 2. Per-locale handback columns: Latest Target File (I), Latest Handback
#    File (J) and Latest Handback DateTime (K) for both tracked files.
# ---------------------------------------------------------------------------

# zh-cn
$zh.Range("I2").Value = $file1Md
$zh.Range("I2").Style = "Hyperlink"
$zh.Range("J2").Value = $file1ZhXlf
$zh.Range("K2").Value = $zhHandbackTime

$zh.Range("I3").Value = $file2Md
$zh.Range("I3").Style = "Hyperlink"
$zh.Range("J3").Value = $file2ZhXlf
$zh.Range("K3").Value = $zhHandbackTime

$zh.Hyperlinks.Add($zh.Range("I2"), $file1Url, "", "", $file1Md)
$zh.Hyperlinks.Add($zh.Range("I3"), $file2Url, "", "", $file2Md)

# de-de
$de.Range("I2").Value = $file1Md
$de.Range("I2").Style = "Hyperlink"
$de.Range("J2").Value = $file1DeXlf
$de.Range("K2").Value = $deHandbackTime

$de.Range("I3").Value = $file2Md
$de.Range("I3").Style = "Hyperlink"
$de.Range("J3").Value = $file2DeXlf
$de.Range("K3").Value = $deHandbackTime

$de.Hyperlinks.Add($de.Range("I2"), $file1Url, "", "", $file1Md)
$de.Hyperlinks.Add($de.Range("I3"), $file2Url, "", "", $file2Md)

# ---------------------------------------------------------------------------
# 3. Widen the columns that now hold longer text so everything stays
#    readable (mirrors the width bump seen after Excel's own autofit).
# ---------------------------------------------------------------------------
$ovw.Columns(5).ColumnWidth = 29.166666666666664
$ovw.Columns(6).ColumnWidth = 29.166666666666664

$zh.Columns(3).ColumnWidth  = 29.166666666666664
$zh.Columns(9).ColumnWidth  = 39.166666666666664
$zh.Columns(10).ColumnWidth = 39.166666666666664

$de.Columns(3).ColumnWidth  = 29.166666666666664
$de.Columns(9).ColumnWidth  = 39.166666666666664
$de.Columns(10).ColumnWidth = 39.166666666666664
